$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.832.16'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '3.524.91'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''598.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").Value = '''143.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("D7").Value = '3.522.86'
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '''0.498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("E10").Value = '  -1.17%  '
$ws.Range("E11").Value = '  -2.31%  '
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '4.126.12'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("E14").Value = '  -3.59%  '
$ws.Range("E15").Value = '  -4.62%  '
$ws.Range("D16").Value = '3.520.82'
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '65.852.24'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").Value = '''10.91'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.37%  '
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '''14.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.81%  '
$ws.Range("D22").Value = '''413.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.08%  '
$ws.Range("D23").Value = '''0.598'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.69%  '
$ws.Range("D24").Value = '''77.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.78%  '
$ws.Range("D25").Value = '3.667.84'
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''7.80'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("D30").Value = '''8.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").Value = '3.523.85'
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("D33").Value = '''0.153'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '''7.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.71%  '
$ws.Range("D37").Value = '''1.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.27%  '
$ws.Range("D38").Value = '''175.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.22%  '
$ws.Range("D39").Value = '''5.27'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.93%  '
$ws.Range("E40").Value = '  -8.13%  '
$ws.Range("D41").Value = '''0.0820'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.10%  '
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '''0.857'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.32%  '
$ws.Range("D44").Value = '''45.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("E45").Value = '  -7.62%  '
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").Value = '''2.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.02%  '
$ws.Range("D48").Value = '''7.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("D49").Value = '''22.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("D51").Value = '''23.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.99%  '
